# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (D5) and
# "Correspond Handback DateTime" (G5) timestamps on the
# zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-28 04:31:58"
$wsZhCn.Range("G5").Value = "2016-01-28 04:32:36"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-28 04:32:08"
$wsDeDe.Range("G5").Value = "2016-01-28 04:32:54"
